# "consistent format for input data"
#
# The import sheet had a stray one-cell header row ("labels" in A1, B1
# empty) and a duplicated data row ("Productive Daily Life" / "V", which
# already existed as "Productive Daily Life" / "C"). Both rows are
# removed so the sheet is a clean, consistent name/category table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 ("labels") - not part of the actual name/ACV data, drop it.
$ws.Rows(1).Delete() | Out-Null

# After the first delete, the duplicate "Productive Daily Life" row
# (originally row 16) has shifted up to row 15 - drop it too.
$ws.Rows(15).Delete() | Out-Null

# Update the window scroll position / selection to reflect where the
# edit left off.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 5
$aw.ScrollColumn = 1
$ws.Rows(15).Select() | Out-Null
